$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh cryptocurrency price/volume snapshot (GitHub Actions scrape update)

$ws.Range("D2").Value = "'67.893.45"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.87%  "

$ws.Range("D3").Value = "'3.817.05"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.07%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'599.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").Value = "'169.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "

$ws.Range("D7").Value = "'3.816.41"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.07%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +0.27%  "

$ws.Range("E10").Value = "  -1.02%  "

$ws.Range("D11").Value = "'6.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.19%  "

$ws.Range("D12").Value = "'0.462"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.73%  "

$ws.Range("D13").Value = "'0.0000276"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.08%  "

$ws.Range("D14").Value = "'37.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "

$ws.Range("D15").Value = "'4.457.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.07%  "

$ws.Range("D16").Value = "'3.808.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.37%  "

$ws.Range("D17").Value = "'18.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.21%  "

$ws.Range("D18").Value = "'67.986.25"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").Value = "'7.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "

$ws.Range("E20").Value = "  +0.15%  "

$ws.Range("D21").Value = "'10.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.25%  "

$ws.Range("D22").Value = "'469.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").Value = "'0.743"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.54%  "

$ws.Range("E24").Value = "  -9.40%  "

$ws.Range("D25").Value = "'83.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.20%  "

$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.65%  "

$ws.Range("D27").Value = "'12.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.27%  "

$ws.Range("D28").Value = "'10.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.56%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  -1.57%  "

$ws.Range("D31").Value = "'3.965.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.04%  "

$ws.Range("E32").Value = "  -1.44%  "

$ws.Range("E33").Value = "  -1.39%  "

$ws.Range("D34").Value = "'30.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.30%  "

$ws.Range("D35").Value = "'9.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").Value = "'3.781.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.27%  "

$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.48%  "

$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.107"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.60%  "

$ws.Range("D39").Value = "'6.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.23%  "

$ws.Range("E40").Value = "  -1.27%  "

$ws.Range("E41").Value = "  -2.31%  "

$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("D43").Value = "'0.320"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.91%  "

$ws.Range("D44").Value = "'8.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.33%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("E46").Value = "  -1.28%  "

$ws.Range("D47").Value = "'411.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.20%  "

$ws.Range("D48").Value = "'46.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D49").Value = "'0.000288"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.37%  "

$ws.Range("D50").Value = "'142.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.71%  "

$ws.Range("E51").Value = "  +0.05%  "
